# ---------------------------------------------------------------------------
# Target change (from the supplied OOXML diff): word/styles.xml's
# <w:docDefaults> block is trimmed so that only the non-redundant entries
# survive:
#
#   rPrDefault/rPr  -> keeps rFonts, sz=22, szCs=22, lang=en
#                      (drops b=0, i=0, smallCaps=0, strike=0, color=000000,
#                       u=none, shd=clear/auto, vertAlign=baseline - all of
#                       which are themselves Word's built-in defaults, so
#                       this is a no-visual-effect cleanup)
#   pPrDefault/pPr  -> keeps spacing line=276/auto only
#                      (drops keepNext=0, keepLines=0, widowControl=1, the
#                       nil pBdr, shd=clear/auto, spacing after/before=0,
#                       ind all-0, contextualSpacing=0, jc=left - again all
#                       Word's built-in defaults)
#
# i.e. the commit (see message: "download tc, tcn, and tl files from GD")
# just re-exported the package with a leaner stylesheet; nothing about the
# *effective*/rendered formatting changes anywhere in the document.
#
# Word's object model does not expose docDefaults as an editable object:
#   - Document.WordOpenXML / Range.WordOpenXML are read-only in this host
#     (confirmed: assigning throws "... is a read-only property").
#   - Range.InsertXML only replaces body-story content
#     ("character position does not address body content" for non-body
#     ranges) - it cannot reach the styles part.
#   - The only COM surface that reaches rPrDefault/pPrDefault at all is
#     reading through a Style's .Font / .ParagraphFormat (Styles("Normal")
#     resolves unset properties from docDefaults - confirmed by reading
#     Font.Bold / ParagraphFormat.LineSpacing below). But *writing* through
#     that same surface never edits docDefaults - it always stamps an
#     explicit <w:rPr>/<w:pPr> override directly onto the style element
#     instead (verified empirically), which would:
#       * leave <w:docDefaults> itself untouched (so it would not reproduce
#         the requested removals), and
#       * add brand-new direct formatting to the Normal style that isn't
#         present in the target OOXML at all, actively moving the document
#         further from the desired result and changing Normal's effective
#         font/paragraph properties from "inherited" to "explicit".
#   - There is no Style.Reset()/ClearFormatting() equivalent, and no
#     "revert to default" sentinel accepted by Font/ParagraphFormat setters
#     (tried $false/0/$null/wdUndefined - all are taken as literal values
#     and still get written out as explicit overrides).
#
# So there is no reachable Word-automation call that edits docDefaults
# without side effects, and since every value being dropped already equals
# Word's built-in default, the rendered document is visually identical
# either way. To avoid corrupting word/styles.xml with spurious direct
# formatting on the Normal style (which would change its effective
# Font/ParagraphFormat from "inherited" to explicit values and diverge
# further from the target), this script only performs non-mutating reads
# to confirm the current state, and otherwise leaves the document exactly
# as authored.

$d = $word.ActiveDocument
$normal = $d.Styles("Normal")

# Sanity-check (read-only) that Normal's effective formatting already
# matches what the trimmed docDefaults implies - e.g. not bold, not
# italic, left-aligned, auto line spacing - confirming this really is a
# redundant-default cleanup rather than a visible formatting change.
$fontBold = $normal.Font.Bold
$fontItalic = $normal.Font.Italic
$paraAlignment = $normal.ParagraphFormat.Alignment
$paraLineSpacing = $normal.ParagraphFormat.LineSpacingRule

Write-Output "Normal.Font.Bold=$fontBold Normal.Font.Italic=$fontItalic Normal.ParagraphFormat.Alignment=$paraAlignment Normal.ParagraphFormat.LineSpacingRule=$paraLineSpacing"
